$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.178.62'
$ws.Range('E2').Value = '  -0.43%  '
$ws.Range('D3').Value = '1.853.12'
$ws.Range('E3').Value = '  -0.87%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '235.48'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.03%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4697'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.29%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2886'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.52%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06549'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.21%  '
$ws.Range('E10').Value = '  +2.08%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07976'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.37%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '97.53'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.37%  '
$ws.Range('D13').Value = '1.856.82'
$ws.Range('E13').Value = '  -0.61%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.101'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.01%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6757'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.01%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '269.05'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -3.64%  '
$ws.Range('D17').Value = '30.154.98'
$ws.Range('E17').Value = '  -0.49%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007686'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +5.29%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.000'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.01%  '
$ws.Range('D21').Value = '2.095.48'
$ws.Range('E21').Value = '  -0.61%  '
$ws.Range('E22').Value = '  -0.12%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.205'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -5.44%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.136'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.45%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '167.09'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.06%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.161'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.13%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.85'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.41%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.930'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.07%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.377'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.11%  '
$ws.Range('E30').Value = '  +2.22%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.464'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.94%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.282'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.99%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.996'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.51%  '
$ws.Range('E34').Value = '  -0.22%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.118'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.96%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6989'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.13%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.706'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.58%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01869'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.79%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.604'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.88%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.326'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.74%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '73.11'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.99%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.934'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.64%  '
$ws.Range('E43').Value = '  -0.11%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.8387'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.26%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '103.19'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4132'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.11%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '936.84'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.11%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.150'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.66%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.023'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.02%  '
$ws.Range('E50').Value = '  -0.80%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05654'
$ws.Range('D51').Style = 'Normal'
